# Fruta / hortaliza, semanal
# A new weekly record is inserted at row 19 of Sheet1, pushing all
# subsequent rows (old rows 19..137) down by one (new rows 20..138).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 19 - this shifts rows 19-137 down to 20-138
$ws.Rows("19").Insert()

# Populate the new row 19 with the new weekly record
$ws.Range("A19").Value2 = 4
$ws.Range("B19").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C19").Value2 = "Los Lagos"
$ws.Range("D19").Value2 = 44462
$ws.Range("E19").Value2 = 10
$ws.Range("F19").Value2 = "Fruta"
$ws.Range("G19").Value2 = 100108
$ws.Range("H19").Value2 = "Tropicales y subtropicales"
$ws.Range("I19").Value2 = 100108005
$ws.Range("J19").Value2 = "Piña"
$ws.Range("K19").Value2 = "Caramelo"
$ws.Range("L19").Value2 = "Primera"
$ws.Range("M19").Value2 = 40
$ws.Range("N19").Value2 = 22000
$ws.Range("O19").Value2 = 22000
$ws.Range("P19").Value2 = 22000
$ws.Range("Q19").Value2 = "$/caja 12 unidades"
$ws.Range("R19").Value2 = "Ecuador"
$ws.Range("S19").Value2 = 1833
$ws.Range("T19").Value2 = 12
